$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D, shifting D:K -> E:L.
$ws.Columns("D").Insert()

# Copy number formats/styles from the (now-shifted) column E into the new
# blank column D so the new year's data renders with the same date /
# number formatting as its neighbours.
$ws.Range("E5:E102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new column D with the latest (2018) financial figures.
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 27700
$ws.Range("D9").Value = "NA"
$ws.Range("D10").Value = "NA"
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("D17").Value = 1200
$ws.Range("D18").Value = 26400
$ws.Range("D20").Value = -13800
$ws.Range("D21").Value = 13500
$ws.Range("D22").Value = 0
$ws.Range("D23").Value = 12600
$ws.Range("D24").Value = 2700
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 10000
$ws.Range("D27").Value = 10000
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 100
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 13800
$ws.Range("D33").Value = 10100
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 10100
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 145800
$ws.Range("D42").Value = 147400
$ws.Range("D43").Value = 0
$ws.Range("D44").Value = 0
$ws.Range("D45").Value = 0
$ws.Range("D46").Value = 0
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 14700
$ws.Range("D49").Value = 0
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 3500
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 746900
$ws.Range("D57").Value = 0
$ws.Range("D58").Value = 0
$ws.Range("D59").Value = 12400
$ws.Range("D60").Value = 0
$ws.Range("D61").Value = 4000
$ws.Range("D62").Value = "NA"
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 669800
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 6300
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 77100
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 10100
$ws.Range("D83").Value = 900
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 12700
$ws.Range("D91").Value = -600
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = 17700
$ws.Range("D96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 46300
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = 76700
